# Trade #202 closed at 2026-02-17 10:08:43 - unknown UNKNOWN +0.000%
#
# - Updates the Summary sheet roll-up numbers
# - Updates the volatility_scorer row on the Strategy Status sheet
# - Appends the new closed trade (#202, volatility_scorer) and the new
#   open trade (#203, MarketMaking) to the "All Trades" log, and mirrors
#   each into its own per-strategy sheet.

$wb = $excel.ActiveWorkbook

function Set-TextCell($range, $val) {
    # Writing a date-/time-looking literal straight into .Value lets Excel's
    # autodetection turn it into a serial-number date cell. Force the cell
    # to Text first, write the literal, then drop back to the default style
    # so we don't leave stray per-cell number formatting behind.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.53
$summary.Range("B4").Value = -0.46
$summary.Range("B5").Value = -0.05
$summary.Range("B6").Value = 202
$summary.Range("B8").Value = 86
$summary.Range("B9").Value = 41.58

# ---------------------------------------------------------------------
# Strategy Status sheet - volatility_scorer row (row 12)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C12").Value = 99.17
$status.Range("D12").Value = 19
$status.Range("E12").Value = -0.83
$status.Range("F12").Value = -0.83
$status.Range("G12").Value = 26.32

# ---------------------------------------------------------------------
# All Trades sheet - append trade #202 (row 203) and #203 (row 204)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$r = 203
$allTrades.Cells.Item($r, 1).Value = 202
Set-TextCell $allTrades.Cells.Item($r, 2) "2026-02-17"
Set-TextCell $allTrades.Cells.Item($r, 3) "10:08:36"
$allTrades.Cells.Item($r, 4).Value = "volatility_scorer"
$allTrades.Cells.Item($r, 5).Value = "NEUTRAL"
$allTrades.Cells.Item($r, 6).Value = 0.02
$allTrades.Cells.Item($r, 7).Value = 0.01
$allTrades.Cells.Item($r, 8).Value = "CLOSED"
$allTrades.Cells.Item($r, 9).Value = -50
$allTrades.Cells.Item($r, 10).Value = -0.01
$allTrades.Cells.Item($r, 11).Value = 99.17
$allTrades.Cells.Item($r, 12).Value = 0
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0.85
$allTrades.Cells.Item($r, 15).Value = "Low vol market (score: inf) - ideal for market making"
$allTrades.Cells.Item($r, 16).Value = "early_exit"
$allTrades.Cells.Item($r, 17).Value = 0.18

$r = 204
$allTrades.Cells.Item($r, 1).Value = 203
Set-TextCell $allTrades.Cells.Item($r, 2) "2026-02-17"
Set-TextCell $allTrades.Cells.Item($r, 3) "10:08:36"
$allTrades.Cells.Item($r, 4).Value = "MarketMaking"
$allTrades.Cells.Item($r, 5).Value = "UP"
$allTrades.Cells.Item($r, 6).Value = 0.98
$allTrades.Cells.Item($r, 7).Value = ""     # no exit price yet - trade still open
$allTrades.Cells.Item($r, 8).Value = "OPEN"
$allTrades.Cells.Item($r, 9).Value = 0
$allTrades.Cells.Item($r, 10).Value = 0
$allTrades.Cells.Item($r, 11).Value = 100.3671991854616
$allTrades.Cells.Item($r, 12).Value = 0
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0.6
$allTrades.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"
$allTrades.Cells.Item($r, 16).Value = ""    # no exit reason yet - trade still open
$allTrades.Cells.Item($r, 17).Value = 0

# ---------------------------------------------------------------------
# volatility_scorer sheet - append trade #202 (row 20)
# ---------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")

$r = 20
$volScorer.Cells.Item($r, 1).Value = 202
Set-TextCell $volScorer.Cells.Item($r, 2) "2026-02-17"
Set-TextCell $volScorer.Cells.Item($r, 3) "10:08:36"
$volScorer.Cells.Item($r, 4).Value = "volatility_scorer"
$volScorer.Cells.Item($r, 5).Value = "NEUTRAL"
$volScorer.Cells.Item($r, 6).Value = 0.02
$volScorer.Cells.Item($r, 7).Value = 0.01
$volScorer.Cells.Item($r, 8).Value = "CLOSED"
$volScorer.Cells.Item($r, 9).Value = -50
$volScorer.Cells.Item($r, 10).Value = -0.01
$volScorer.Cells.Item($r, 11).Value = 99.17
$volScorer.Cells.Item($r, 12).Value = 0
$volScorer.Cells.Item($r, 13).Value = 0
$volScorer.Cells.Item($r, 14).Value = 0.85
$volScorer.Cells.Item($r, 15).Value = "Low vol market (score: inf) - ideal for market making"
$volScorer.Cells.Item($r, 16).Value = "early_exit"
$volScorer.Cells.Item($r, 17).Value = 0.18

# ---------------------------------------------------------------------
# MarketMaking sheet - append trade #203 (row 185)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$r = 185
$marketMaking.Cells.Item($r, 1).Value = 203
Set-TextCell $marketMaking.Cells.Item($r, 2) "2026-02-17"
Set-TextCell $marketMaking.Cells.Item($r, 3) "10:08:36"
$marketMaking.Cells.Item($r, 4).Value = "MarketMaking"
$marketMaking.Cells.Item($r, 5).Value = "UP"
$marketMaking.Cells.Item($r, 6).Value = 0.98
$marketMaking.Cells.Item($r, 7).Value = ""  # no exit price yet - trade still open
$marketMaking.Cells.Item($r, 8).Value = "OPEN"
$marketMaking.Cells.Item($r, 9).Value = 0
$marketMaking.Cells.Item($r, 10).Value = 0
$marketMaking.Cells.Item($r, 11).Value = 100.3671991854616
$marketMaking.Cells.Item($r, 12).Value = 0
$marketMaking.Cells.Item($r, 13).Value = 0
$marketMaking.Cells.Item($r, 14).Value = 0.6
$marketMaking.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item($r, 16).Value = "" # no exit reason yet - trade still open
$marketMaking.Cells.Item($r, 17).Value = 0
